$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Mdk"
$ws.Cells.Item(2,3).Value = "Itga4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = [double]"3"
$ws.Cells.Item(2,6).Value = [double]"1"
$ws.Cells.Item(2,7).Value = [double]"1.324023666666666"
$ws.Cells.Item(2,8).Value = [double]"3.972071"
$ws.Cells.Item(2,9).Value = [double]"0.01518042398701374"
$ws.Cells.Item(2,10).Value = [double]"0.01518042398701374"
$ws.Cells.Item(2,11).Value = [double]"2"
$ws.Cells.Item(2,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2,13).Value = [double]"0.3809099999999999"
$ws.Cells.Item(2,14).Value = [double]"1.14273"
$ws.Cells.Item(2,15).Value = [double]"0.006635732896411959"
$ws.Cells.Item(2,16).Value = [double]"0.006635732896411961"
$ws.Cells.Item(2,17).Value = [double]"0.5043338548699998"
$ws.Cells.Item(2,18).Value = [double]"4.539004693829999"
$ws.Cells.Item(2,19).Value = [double]"0.0001007332388321083"
$ws.Cells.Item(2,20).Value = [double]"0.0001007332388321083"

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Mdk"
$ws.Cells.Item(3,3).Value = "Itga4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = [double]"3"
$ws.Cells.Item(3,6).Value = [double]"1"
$ws.Cells.Item(3,7).Value = [double]"1.324023666666666"
$ws.Cells.Item(3,8).Value = [double]"3.972071"
$ws.Cells.Item(3,9).Value = [double]"0.01518042398701374"
$ws.Cells.Item(3,10).Value = [double]"0.01518042398701374"
$ws.Cells.Item(3,11).Value = [double]"3"
$ws.Cells.Item(3,12).Value = [double]"1"
$ws.Cells.Item(3,13).Value = [double]"0.3194813333333333"
$ws.Cells.Item(3,14).Value = [double]"0.9584440000000001"
$ws.Cells.Item(3,15).Value = [double]"0.005565600255676026"
$ws.Cells.Item(3,16).Value = [double]"0.005565600255676028"
$ws.Cells.Item(3,17).Value = [double]"0.4230008463915555"
$ws.Cells.Item(3,18).Value = [double]"3.807007617524"
$ws.Cells.Item(3,19).Value = [double]"8.448817162339413E-05"
$ws.Cells.Item(3,20).Value = [double]"8.448817162339416E-05"

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Mdk"
$ws.Cells.Item(4,3).Value = "Itga4"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = [double]"3"
$ws.Cells.Item(4,6).Value = [double]"1"
$ws.Cells.Item(4,7).Value = [double]"1.324023666666666"
$ws.Cells.Item(4,8).Value = [double]"3.972071"
$ws.Cells.Item(4,9).Value = [double]"0.01518042398701374"
$ws.Cells.Item(4,10).Value = [double]"0.01518042398701374"
$ws.Cells.Item(4,11).Value = [double]"3"
$ws.Cells.Item(4,12).Value = [double]"1"
$ws.Cells.Item(4,13).Value = [double]"54.69403966666666"
$ws.Cells.Item(4,14).Value = [double]"164.082119"
$ws.Cells.Item(4,15).Value = [double]"0.9528104755815301"
$ws.Cells.Item(4,16).Value = [double]"0.9528104755815303"
$ws.Cells.Item(4,17).Value = [double]"72.4162029442721"
$ws.Cells.Item(4,18).Value = [double]"651.7458264984489"
$ws.Cells.Item(4,19).Value = [double]"0.01446406699859583"
$ws.Cells.Item(4,20).Value = [double]"0.01446406699859583"

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Mdk"
$ws.Cells.Item(5,3).Value = "Itga4"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = [double]"3"
$ws.Cells.Item(5,6).Value = [double]"1"
$ws.Cells.Item(5,7).Value = [double]"1.324023666666666"
$ws.Cells.Item(5,8).Value = [double]"3.972071"
$ws.Cells.Item(5,9).Value = [double]"0.01518042398701374"
$ws.Cells.Item(5,10).Value = [double]"0.01518042398701374"
$ws.Cells.Item(5,11).Value = [double]"3"
$ws.Cells.Item(5,12).Value = [double]"1"
$ws.Cells.Item(5,13).Value = [double]"2.008422"
$ws.Cells.Item(5,14).Value = [double]"6.025265999999999"
$ws.Cells.Item(5,15).Value = [double]"0.03498819126638183"
$ws.Cells.Item(5,16).Value = [double]"0.03498819126638184"
$ws.Cells.Item(5,17).Value = [double]"2.659198260654"
$ws.Cells.Item(5,18).Value = [double]"23.932784345886"
$ws.Cells.Item(5,19).Value = [double]"0.0005311355779624073"
$ws.Cells.Item(5,20).Value = [double]"0.0005311355779624074"

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Mdk"
$ws.Cells.Item(6,3).Value = "Itga4"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = [double]"3"
$ws.Cells.Item(6,6).Value = [double]"1"
$ws.Cells.Item(6,7).Value = [double]"81.17653533333333"
$ws.Cells.Item(6,8).Value = [double]"243.529606"
$ws.Cells.Item(6,9).Value = [double]"0.9307191821270077"
$ws.Cells.Item(6,10).Value = [double]"0.9307191821270075"
$ws.Cells.Item(6,11).Value = [double]"2"
$ws.Cells.Item(6,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(6,13).Value = [double]"0.3809099999999999"
$ws.Cells.Item(6,14).Value = [double]"1.14273"
$ws.Cells.Item(6,15).Value = [double]"0.006635732896411959"
$ws.Cells.Item(6,16).Value = [double]"0.006635732896411961"
$ws.Cells.Item(6,17).Value = [double]"30.92095407381999"
$ws.Cells.Item(6,18).Value = [double]"278.28858666438"
$ws.Cells.Item(6,19).Value = [double]"0.006176003894161818"
$ws.Cells.Item(6,20).Value = [double]"0.006176003894161819"

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Mdk"
$ws.Cells.Item(7,3).Value = "Itga4"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = [double]"3"
$ws.Cells.Item(7,6).Value = [double]"1"
$ws.Cells.Item(7,7).Value = [double]"81.17653533333333"
$ws.Cells.Item(7,8).Value = [double]"243.529606"
$ws.Cells.Item(7,9).Value = [double]"0.9307191821270077"
$ws.Cells.Item(7,10).Value = [double]"0.9307191821270075"
$ws.Cells.Item(7,11).Value = [double]"3"
$ws.Cells.Item(7,12).Value = [double]"1"
$ws.Cells.Item(7,13).Value = [double]"0.3194813333333333"
$ws.Cells.Item(7,14).Value = [double]"0.9584440000000001"
$ws.Cells.Item(7,15).Value = [double]"0.005565600255676026"
$ws.Cells.Item(7,16).Value = [double]"0.005565600255676028"
$ws.Cells.Item(7,17).Value = [double]"25.93438774367378"
$ws.Cells.Item(7,18).Value = [double]"233.409489693064"
$ws.Cells.Item(7,19).Value = [double]"0.005180010918008656"
$ws.Cells.Item(7,20).Value = [double]"0.005180010918008657"

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Mdk"
$ws.Cells.Item(8,3).Value = "Itga4"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = [double]"3"
$ws.Cells.Item(8,6).Value = [double]"1"
$ws.Cells.Item(8,7).Value = [double]"81.17653533333333"
$ws.Cells.Item(8,8).Value = [double]"243.529606"
$ws.Cells.Item(8,9).Value = [double]"0.9307191821270077"
$ws.Cells.Item(8,10).Value = [double]"0.9307191821270075"
$ws.Cells.Item(8,11).Value = [double]"3"
$ws.Cells.Item(8,12).Value = [double]"1"
$ws.Cells.Item(8,13).Value = [double]"54.69403966666666"
$ws.Cells.Item(8,14).Value = [double]"164.082119"
$ws.Cells.Item(8,15).Value = [double]"0.9528104755815301"
$ws.Cells.Item(8,16).Value = [double]"0.9528104755815303"
$ws.Cells.Item(8,17).Value = [double]"4439.872643523901"
$ws.Cells.Item(8,18).Value = [double]"39958.85379171511"
$ws.Cells.Item(8,19).Value = [double]"0.8867989865552869"
$ws.Cells.Item(8,20).Value = [double]"0.886798986555287"

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Mdk"
$ws.Cells.Item(9,3).Value = "Itga4"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = [double]"3"
$ws.Cells.Item(9,6).Value = [double]"1"
$ws.Cells.Item(9,7).Value = [double]"81.17653533333333"
$ws.Cells.Item(9,8).Value = [double]"243.529606"
$ws.Cells.Item(9,9).Value = [double]"0.9307191821270077"
$ws.Cells.Item(9,10).Value = [double]"0.9307191821270075"
$ws.Cells.Item(9,11).Value = [double]"3"
$ws.Cells.Item(9,12).Value = [double]"1"
$ws.Cells.Item(9,13).Value = [double]"2.008422"
$ws.Cells.Item(9,14).Value = [double]"6.025265999999999"
$ws.Cells.Item(9,15).Value = [double]"0.03498819126638183"
$ws.Cells.Item(9,16).Value = [double]"0.03498819126638184"
$ws.Cells.Item(9,17).Value = [double]"163.036739447244"
$ws.Cells.Item(9,18).Value = [double]"1467.330655025196"
$ws.Cells.Item(9,19).Value = [double]"0.03256418075955021"
$ws.Cells.Item(9,20).Value = [double]"0.03256418075955021"

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Mdk"
$ws.Cells.Item(10,3).Value = "Itga4"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = [double]"3"
$ws.Cells.Item(10,6).Value = [double]"1"
$ws.Cells.Item(10,7).Value = [double]"1.192675"
$ws.Cells.Item(10,8).Value = [double]"3.578025"
$ws.Cells.Item(10,9).Value = [double]"0.0136744626508778"
$ws.Cells.Item(10,10).Value = [double]"0.0136744626508778"
$ws.Cells.Item(10,11).Value = [double]"2"
$ws.Cells.Item(10,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10,13).Value = [double]"0.3809099999999999"
$ws.Cells.Item(10,14).Value = [double]"1.14273"
$ws.Cells.Item(10,15).Value = [double]"0.006635732896411959"
$ws.Cells.Item(10,16).Value = [double]"0.006635732896411961"
$ws.Cells.Item(10,17).Value = [double]"0.45430183425"
$ws.Cells.Item(10,18).Value = [double]"4.088716508249999"
$ws.Cells.Item(10,19).Value = [double]"9.074008165318652E-05"
$ws.Cells.Item(10,20).Value = [double]"9.074008165318653E-05"

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Mdk"
$ws.Cells.Item(11,3).Value = "Itga4"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = [double]"3"
$ws.Cells.Item(11,6).Value = [double]"1"
$ws.Cells.Item(11,7).Value = [double]"1.192675"
$ws.Cells.Item(11,8).Value = [double]"3.578025"
$ws.Cells.Item(11,9).Value = [double]"0.0136744626508778"
$ws.Cells.Item(11,10).Value = [double]"0.0136744626508778"
$ws.Cells.Item(11,11).Value = [double]"3"
$ws.Cells.Item(11,12).Value = [double]"1"
$ws.Cells.Item(11,13).Value = [double]"0.3194813333333333"
$ws.Cells.Item(11,14).Value = [double]"0.9584440000000001"
$ws.Cells.Item(11,15).Value = [double]"0.005565600255676026"
$ws.Cells.Item(11,16).Value = [double]"0.005565600255676028"
$ws.Cells.Item(11,17).Value = [double]"0.3810373992333334"
$ws.Cells.Item(11,18).Value = [double]"3.4293365931"
$ws.Cells.Item(11,19).Value = [double]"7.610659282595776E-05"
$ws.Cells.Item(11,20).Value = [double]"7.610659282595778E-05"

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Mdk"
$ws.Cells.Item(12,3).Value = "Itga4"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = [double]"3"
$ws.Cells.Item(12,6).Value = [double]"1"
$ws.Cells.Item(12,7).Value = [double]"1.192675"
$ws.Cells.Item(12,8).Value = [double]"3.578025"
$ws.Cells.Item(12,9).Value = [double]"0.0136744626508778"
$ws.Cells.Item(12,10).Value = [double]"0.0136744626508778"
$ws.Cells.Item(12,11).Value = [double]"3"
$ws.Cells.Item(12,12).Value = [double]"1"
$ws.Cells.Item(12,13).Value = [double]"54.69403966666666"
$ws.Cells.Item(12,14).Value = [double]"164.082119"
$ws.Cells.Item(12,15).Value = [double]"0.9528104755815301"
$ws.Cells.Item(12,16).Value = [double]"0.9528104755815303"
$ws.Cells.Item(12,17).Value = [double]"65.23221375944166"
$ws.Cells.Item(12,18).Value = [double]"587.0899238349749"
$ws.Cells.Item(12,19).Value = [double]"0.01302917126170475"
$ws.Cells.Item(12,20).Value = [double]"0.01302917126170475"

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Mdk"
$ws.Cells.Item(13,3).Value = "Itga4"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = [double]"3"
$ws.Cells.Item(13,6).Value = [double]"1"
$ws.Cells.Item(13,7).Value = [double]"1.192675"
$ws.Cells.Item(13,8).Value = [double]"3.578025"
$ws.Cells.Item(13,9).Value = [double]"0.0136744626508778"
$ws.Cells.Item(13,10).Value = [double]"0.0136744626508778"
$ws.Cells.Item(13,11).Value = [double]"3"
$ws.Cells.Item(13,12).Value = [double]"1"
$ws.Cells.Item(13,13).Value = [double]"2.008422"
$ws.Cells.Item(13,14).Value = [double]"6.025265999999999"
$ws.Cells.Item(13,15).Value = [double]"0.03498819126638183"
$ws.Cells.Item(13,16).Value = [double]"0.03498819126638184"
$ws.Cells.Item(13,17).Value = [double]"2.39539470885"
$ws.Cells.Item(13,18).Value = [double]"21.55855237965"
$ws.Cells.Item(13,19).Value = [double]"0.0004784447146939073"
$ws.Cells.Item(13,20).Value = [double]"0.0004784447146939073"

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Mdk"
$ws.Cells.Item(14,3).Value = "Itga4"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = [double]"3"
$ws.Cells.Item(14,6).Value = [double]"1"
$ws.Cells.Item(14,7).Value = [double]"3.525915333333334"
$ws.Cells.Item(14,8).Value = [double]"10.577746"
$ws.Cells.Item(14,9).Value = [double]"0.04042593123510095"
$ws.Cells.Item(14,10).Value = [double]"0.04042593123510094"
$ws.Cells.Item(14,11).Value = [double]"2"
$ws.Cells.Item(14,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(14,13).Value = [double]"0.3809099999999999"
$ws.Cells.Item(14,14).Value = [double]"1.14273"
$ws.Cells.Item(14,15).Value = [double]"0.006635732896411959"
$ws.Cells.Item(14,16).Value = [double]"0.006635732896411961"
$ws.Cells.Item(14,17).Value = [double]"1.34305640962"
$ws.Cells.Item(14,18).Value = [double]"12.08750768658"
$ws.Cells.Item(14,19).Value = [double]"0.0002682556817648471"
$ws.Cells.Item(14,20).Value = [double]"0.0002682556817648472"

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Mdk"
$ws.Cells.Item(15,3).Value = "Itga4"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = [double]"3"
$ws.Cells.Item(15,6).Value = [double]"1"
$ws.Cells.Item(15,7).Value = [double]"3.525915333333334"
$ws.Cells.Item(15,8).Value = [double]"10.577746"
$ws.Cells.Item(15,9).Value = [double]"0.04042593123510095"
$ws.Cells.Item(15,10).Value = [double]"0.04042593123510094"
$ws.Cells.Item(15,11).Value = [double]"3"
$ws.Cells.Item(15,12).Value = [double]"1"
$ws.Cells.Item(15,13).Value = [double]"0.3194813333333333"
$ws.Cells.Item(15,14).Value = [double]"0.9584440000000001"
$ws.Cells.Item(15,15).Value = [double]"0.005565600255676026"
$ws.Cells.Item(15,16).Value = [double]"0.005565600255676028"
$ws.Cells.Item(15,17).Value = [double]"1.126464131913778"
$ws.Cells.Item(15,18).Value = [double]"10.138177187224"
$ws.Cells.Item(15,19).Value = [double]"0.0002249945732180193"
$ws.Cells.Item(15,20).Value = [double]"0.0002249945732180193"

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Mdk"
$ws.Cells.Item(16,3).Value = "Itga4"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = [double]"3"
$ws.Cells.Item(16,6).Value = [double]"1"
$ws.Cells.Item(16,7).Value = [double]"3.525915333333334"
$ws.Cells.Item(16,8).Value = [double]"10.577746"
$ws.Cells.Item(16,9).Value = [double]"0.04042593123510095"
$ws.Cells.Item(16,10).Value = [double]"0.04042593123510094"
$ws.Cells.Item(16,11).Value = [double]"3"
$ws.Cells.Item(16,12).Value = [double]"1"
$ws.Cells.Item(16,13).Value = [double]"54.69403966666666"
$ws.Cells.Item(16,14).Value = [double]"164.082119"
$ws.Cells.Item(16,15).Value = [double]"0.9528104755815301"
$ws.Cells.Item(16,16).Value = [double]"0.9528104755815303"
$ws.Cells.Item(16,17).Value = [double]"192.8465531026416"
$ws.Cells.Item(16,18).Value = [double]"1735.618977923774"
$ws.Cells.Item(16,19).Value = [double]"0.03851825076594277"
$ws.Cells.Item(16,20).Value = [double]"0.03851825076594277"

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Mdk"
$ws.Cells.Item(17,3).Value = "Itga4"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = [double]"3"
$ws.Cells.Item(17,6).Value = [double]"1"
$ws.Cells.Item(17,7).Value = [double]"3.525915333333334"
$ws.Cells.Item(17,8).Value = [double]"10.577746"
$ws.Cells.Item(17,9).Value = [double]"0.04042593123510095"
$ws.Cells.Item(17,10).Value = [double]"0.04042593123510094"
$ws.Cells.Item(17,11).Value = [double]"3"
$ws.Cells.Item(17,12).Value = [double]"1"
$ws.Cells.Item(17,13).Value = [double]"2.008422"
$ws.Cells.Item(17,14).Value = [double]"6.025265999999999"
$ws.Cells.Item(17,15).Value = [double]"0.03498819126638183"
$ws.Cells.Item(17,16).Value = [double]"0.03498819126638184"
$ws.Cells.Item(17,17).Value = [double]"7.081525925604001"
$ws.Cells.Item(17,18).Value = [double]"63.733733330436"
$ws.Cells.Item(17,19).Value = [double]"0.001414430214175312"
$ws.Cells.Item(17,20).Value = [double]"0.001414430214175312"
